# Auto-generated edit script: updates LeveProfit calculation sheets
# (currentAveragePrice / Price / Profit columns H..N) per scheduled runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 237.9375
$ws.Range("I6").Value = 187.13333
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 561.39999
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = -449.39999
$ws.Range("N6").Value = -3224

$ws.Range("H129").Value = 1032.0182
$ws.Range("I129").Value = 681.2
$ws.Range("J129").Value = 1109.9778
$ws.Range("K129").Value = 2043.6
$ws.Range("L129").Value = 3329.9334
$ws.Range("M129").Value = 2956.4
$ws.Range("N129").Value = -13329.9334

$ws.Range("H137").Value = 1379.9828
$ws.Range("I137").Value = 1293.1063
$ws.Range("J137").Value = 1751.1818
$ws.Range("K137").Value = 3879.3189
$ws.Range("L137").Value = 5253.5454
$ws.Range("M137").Value = -1329.3189
$ws.Range("N137").Value = -10353.5454

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7690.8115
$ws.Range("I32").Value = 4624.0327
$ws.Range("J32").Value = 31075
$ws.Range("K32").Value = 4624.0327
$ws.Range("L32").Value = 31075
$ws.Range("M32").Value = -4337.0327
$ws.Range("N32").Value = -31649

$ws.Range("H45").Value = 5352.409
$ws.Range("I45").Value = 6102.0527
$ws.Range("J45").Value = 604.6667
$ws.Range("K45").Value = 6102.0527
$ws.Range("L45").Value = 604.6667
$ws.Range("M45").Value = -5725.0527
$ws.Range("N45").Value = -1358.6667

$ws.Range("H74").Value = 1753.7391
$ws.Range("I74").Value = 1741.2307
$ws.Range("J74").Value = 1770
$ws.Range("K74").Value = 1741.2307
$ws.Range("L74").Value = 1770
$ws.Range("M74").Value = -867.2307000000001
$ws.Range("N74").Value = -3518

$ws.Range("H77").Value = 1753.7391
$ws.Range("I77").Value = 1741.2307
$ws.Range("J77").Value = 1770
$ws.Range("K77").Value = 8706.1535
$ws.Range("L77").Value = 8850
$ws.Range("M77").Value = -4338.1535
$ws.Range("N77").Value = -17586

$ws.Range("H110").Value = 1075.9333
$ws.Range("I110").Value = 594.9167
$ws.Range("J110").Value = 3000
$ws.Range("K110").Value = 594.9167
$ws.Range("L110").Value = 3000
$ws.Range("M110").Value = 1450.0833
$ws.Range("N110").Value = -7090

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8573.120999999999
$ws.Range("I20").Value = 802.58826
$ws.Range("J20").Value = 16829.312
$ws.Range("K20").Value = 802.58826
$ws.Range("L20").Value = 16829.312
$ws.Range("M20").Value = -555.58826
$ws.Range("N20").Value = -17323.312

$ws.Range("H105").Value = 16478.666
$ws.Range("I105").Value = 67936.664
$ws.Range("J105").Value = 3614.1667
$ws.Range("K105").Value = 67936.664
$ws.Range("L105").Value = 3614.1667
$ws.Range("M105").Value = -66189.664
$ws.Range("N105").Value = -7108.1667

$ws.Range("H107").Value = 830.4194
$ws.Range("I107").Value = 631.4091
$ws.Range("J107").Value = 1316.8889
$ws.Range("K107").Value = 631.4091
$ws.Range("L107").Value = 1316.8889
$ws.Range("M107").Value = 1288.5909
$ws.Range("N107").Value = -5156.8889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1985.2667
$ws.Range("I16").Value = 1716.6
$ws.Range("J16").Value = 2522.6
$ws.Range("K16").Value = 1716.6
$ws.Range("L16").Value = 2522.6
$ws.Range("M16").Value = -1429.6
$ws.Range("N16").Value = -3096.6

$ws.Range("H31").Value = 3901.1724
$ws.Range("I31").Value = 1571.2954
$ws.Range("J31").Value = 11223.643
$ws.Range("K31").Value = 1571.2954
$ws.Range("L31").Value = 11223.643
$ws.Range("M31").Value = -1276.2954
$ws.Range("N31").Value = -11813.643

$ws.Range("H34").Value = 3901.1724
$ws.Range("I34").Value = 1571.2954
$ws.Range("J34").Value = 11223.643
$ws.Range("K34").Value = 1571.2954
$ws.Range("L34").Value = 11223.643
$ws.Range("M34").Value = -1369.2954
$ws.Range("N34").Value = -11627.643

$ws.Range("H105").Value = 1805.619
$ws.Range("I105").Value = 2237
$ws.Range("J105").Value = 1413.4546
$ws.Range("K105").Value = 2237
$ws.Range("L105").Value = 1413.4546
$ws.Range("M105").Value = -490
$ws.Range("N105").Value = -4907.4546

$ws.Range("H113").Value = 1985.2667
$ws.Range("I113").Value = 1716.6
$ws.Range("J113").Value = 2522.6
$ws.Range("K113").Value = 1716.6
$ws.Range("L113").Value = 2522.6
$ws.Range("M113").Value = 453.4000000000001
$ws.Range("N113").Value = -6862.6

$ws.Range("H132").Value = 2873.2
$ws.Range("I132").Value = 3878
$ws.Range("J132").Value = 2507.818
$ws.Range("K132").Value = 11634
$ws.Range("L132").Value = 7523.454000000001
$ws.Range("M132").Value = -9104
$ws.Range("N132").Value = -12583.454

$ws.Range("H134").Value = 3107.4324
$ws.Range("I134").Value = 3074.4062
$ws.Range("J134").Value = 3318.8
$ws.Range("K134").Value = 9223.2186
$ws.Range("L134").Value = 9956.400000000001
$ws.Range("M134").Value = -6688.2186
$ws.Range("N134").Value = -15026.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 130.125
$ws.Range("I7").Value = 106.833336
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 320.500008
$ws.Range("L7").Value = 600
$ws.Range("M7").Value = -208.500008
$ws.Range("N7").Value = -824

$ws.Range("H80").Value = 1549.4
$ws.Range("I80").Value = 1000
$ws.Range("J80").Value = 1686.75
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 5060.25
$ws.Range("M80").Value = -2064
$ws.Range("N80").Value = -6932.25

$ws.Range("H83").Value = 1549.4
$ws.Range("I83").Value = 1000
$ws.Range("J83").Value = 1686.75
$ws.Range("K83").Value = 9000
$ws.Range("L83").Value = 15180.75
$ws.Range("M83").Value = -4320
$ws.Range("N83").Value = -24540.75

$ws.Range("H92").Value = 597.6667
$ws.Range("I92").Value = 800
$ws.Range("J92").Value = 496.5
$ws.Range("K92").Value = 2400
$ws.Range("L92").Value = 1489.5
$ws.Range("M92").Value = -1152
$ws.Range("N92").Value = -3985.5

$ws.Range("H107").Value = 615.5
$ws.Range("I107").Value = 287.25
$ws.Range("J107").Value = 709.2857
$ws.Range("K107").Value = 861.75
$ws.Range("L107").Value = 2127.8571
$ws.Range("M107").Value = 1058.25
$ws.Range("N107").Value = -5967.8571

$ws.Range("H114").Value = 5941.864
$ws.Range("I114").Value = 551.7778
$ws.Range("J114").Value = 9673.462
$ws.Range("K114").Value = 1655.3334
$ws.Range("L114").Value = 29020.386
$ws.Range("M114").Value = 1598.6666
$ws.Range("N114").Value = -35528.386

$ws.Range("H129").Value = 1118
$ws.Range("I129").Value = 999.625
$ws.Range("J129").Value = 1307.4
$ws.Range("K129").Value = 2998.875
$ws.Range("L129").Value = 3922.2
$ws.Range("M129").Value = 2001.125
$ws.Range("N129").Value = -13922.2

$ws.Range("H131").Value = 963.125
$ws.Range("I131").Value = 703.1667
$ws.Range("J131").Value = 1049.7778
$ws.Range("K131").Value = 2109.5001
$ws.Range("L131").Value = 3149.3334
$ws.Range("M131").Value = 2930.4999
$ws.Range("N131").Value = -13229.3334

$ws.Range("H137").Value = 10313.255
$ws.Range("I137").Value = 8803.267
$ws.Range("J137").Value = 11021.0625
$ws.Range("K137").Value = 26409.801
$ws.Range("L137").Value = 33063.1875
$ws.Range("M137").Value = -21309.801
$ws.Range("N137").Value = -43263.1875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 0
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H113").Value = 40001340
$ws.Range("I113").Value = 90910080
$ws.Range("J113").Value = 1615.8572
$ws.Range("K113").Value = 90910080
$ws.Range("L113").Value = 1615.8572
$ws.Range("M113").Value = -90907910
$ws.Range("N113").Value = -5955.8572

$ws.Range("H126").Value = 4555.625
$ws.Range("I126").Value = 9508.538
$ws.Range("J126").Value = 2170.889
$ws.Range("K126").Value = 28525.614
$ws.Range("L126").Value = 6512.667
$ws.Range("M126").Value = -26055.614
$ws.Range("N126").Value = -11452.667

$ws.Range("H132").Value = 2441.4546
$ws.Range("I132").Value = 2010.6
$ws.Range("J132").Value = 2800.5
$ws.Range("K132").Value = 6031.799999999999
$ws.Range("L132").Value = 8401.5
$ws.Range("M132").Value = -3501.799999999999
$ws.Range("N132").Value = -13461.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 55558204
$ws.Range("I40").Value = 66669064
$ws.Range("J40").Value = 3898.3333
$ws.Range("K40").Value = 66669064
$ws.Range("L40").Value = 3898.3333
$ws.Range("M40").Value = -66668928
$ws.Range("N40").Value = -4170.3333

$ws.Range("H68").Value = 55612732
$ws.Range("I68").Value = 78646.30499999999
$ws.Range("J68").Value = 200001360
$ws.Range("K68").Value = 78646.30499999999
$ws.Range("L68").Value = 200001360
$ws.Range("M68").Value = -77897.30499999999
$ws.Range("N68").Value = -200002858

$ws.Range("H71").Value = 55612732
$ws.Range("I71").Value = 78646.30499999999
$ws.Range("J71").Value = 200001360
$ws.Range("K71").Value = 393231.525
$ws.Range("L71").Value = 1000006800
$ws.Range("M71").Value = -389487.525
$ws.Range("N71").Value = -1000014288

$ws.Range("H132").Value = 12146056
$ws.Range("I132").Value = 16193620
$ws.Range("J132").Value = 3363.0908
$ws.Range("K132").Value = 48580860
$ws.Range("L132").Value = 10089.2724
$ws.Range("M132").Value = -48578330
$ws.Range("N132").Value = -15149.2724

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4019.2
$ws.Range("I62").Value = 3774
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 3774
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -3150
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 4019.2
$ws.Range("I65").Value = 3774
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 18870
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -15750
$ws.Range("N65").Value = -31240

$ws.Range("H81").Value = 1673.2727
$ws.Range("I81").Value = 1295.8334
$ws.Range("J81").Value = 2126.2
$ws.Range("K81").Value = 2591.6668
$ws.Range("L81").Value = 4252.4
$ws.Range("M81").Value = -1530.6668
$ws.Range("N81").Value = -6374.4

$ws.Range("H84").Value = 1673.2727
$ws.Range("I84").Value = 1295.8334
$ws.Range("J84").Value = 2126.2
$ws.Range("K84").Value = 12958.334
$ws.Range("L84").Value = 21262
$ws.Range("M84").Value = -7654.333999999999
$ws.Range("N84").Value = -31870

$ws.Range("H96").Value = 2640
$ws.Range("I96").Value = 2800
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 2800
$ws.Range("L96").Value = 2000
$ws.Range("M96").Value = -1427
$ws.Range("N96").Value = -4746

$ws.Range("H113").Value = 773.0714
$ws.Range("I113").Value = 469.17648
$ws.Range("J113").Value = 1242.7273
$ws.Range("K113").Value = 1407.52944
$ws.Range("L113").Value = 3728.1819
$ws.Range("M113").Value = 762.47056
$ws.Range("N113").Value = -8068.1819

$ws.Range("H122").Value = 1062.409
$ws.Range("I122").Value = 1118.3529
$ws.Range("J122").Value = 872.2
$ws.Range("K122").Value = 3355.0587
$ws.Range("L122").Value = 2616.6
$ws.Range("M122").Value = -905.0587000000005
$ws.Range("N122").Value = -7516.6

$ws.Range("H126").Value = 1156.5714
$ws.Range("I126").Value = 719.2
$ws.Range("J126").Value = 2250
$ws.Range("K126").Value = 2157.6
$ws.Range("L126").Value = 6750
$ws.Range("M126").Value = 312.3999999999996
$ws.Range("N126").Value = -11690

$ws.Range("H132").Value = 1324.1282
$ws.Range("I132").Value = 949.55554
$ws.Range("J132").Value = 2166.9167
$ws.Range("K132").Value = 2848.66662
$ws.Range("L132").Value = 6500.750100000001
$ws.Range("M132").Value = -318.66662
$ws.Range("N132").Value = -11560.7501
